$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 23:41"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 7131060
$ws.Range("C4").Value = 33123
$ws.Range("D4").Value = 4381466
$ws.Range("E4").Value = 2543257
$ws.Range("G4").Value = 856
$ws.Range("H4").Value = 206337

# Row 6: Brasil -> Brasil
$ws.Range("B6").Value = 4624885
$ws.Range("C6").Value = 29550
$ws.Range("E6").Value = 540281
$ws.Range("G6").Value = 818
$ws.Range("H6").Value = 138977

# Row 27: Israel -> Israel
$ws.Range("B27").Value = 204690
$ws.Range("C27").Value = 11316
$ws.Range("D27").Value = 144963
$ws.Range("E27").Value = 58402
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = 1325

# Row 56: Barein -> Barein
$ws.Range("B56").Value = 67014
$ws.Range("C56").Value = 612
$ws.Range("D56").Value = 60117
$ws.Range("E56").Value = 6666

# Row 70: Kenia -> Kenia
$ws.Range("B70").Value = 37348
$ws.Range("C70").Value = 130
$ws.Range("D70").Value = 24253
$ws.Range("E70").Value = 12431
$ws.Range("G70").Value = 5
$ws.Range("H70").Value = 664

# Row 82: Camerun -> Camerun
$ws.Range("B82").Value = 20690
$ws.Range("C82").Value = 92
$ws.Range("E82").Value = 1150

# Row 85: Bulgaria -> Bulgaria
$ws.Range("B85").Value = 19283
$ws.Range("C85").Value = 160
$ws.Range("D85").Value = 13867
$ws.Range("E85").Value = 4637
$ws.Range("G85").Value = 12
$ws.Range("H85").Value = 779

# Row 91: Zambia -> Zambia
$ws.Range("B91").Value = 14443
$ws.Range("C91").Value = 54
$ws.Range("E91").Value = 482
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 332

# Row 96: Namibia -> Namibia
$ws.Range("B96").Value = 10663
$ws.Range("C96").Value = 56
$ws.Range("D96").Value = 8431
$ws.Range("E96").Value = 2115
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 117

# Row 99: Guinea -> Guinea
$ws.Range("B99").Value = 10434
$ws.Range("C99").Value = 47
$ws.Range("D99").Value = 9801
$ws.Range("E99").Value = 568

# Row 105: Gabon -> Gabon
$ws.Range("B105").Value = 8716
$ws.Range("C105").Value = 12
$ws.Range("D105").Value = 7906
$ws.Range("E105").Value = 756

# Row 108: Zimbabue -> Zimbabue
$ws.Range("B108").Value = 7725
$ws.Range("C108").Value = 14
$ws.Range("D108").Value = 6007
$ws.Range("E108").Value = 1491
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 227

# Row 119: Suazilandia -> Suazilandia
$ws.Range("B119").Value = 5343
$ws.Range("C119").Value = 36
$ws.Range("D119").Value = 4693
$ws.Range("E119").Value = 542
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 108

# Row 124: Congo -> Congo
$ws.Range("B124").Value = 5005
$ws.Range("C124").Value = 3
$ws.Range("E124").Value = 1029

# Row 134: Siria -> Siria
$ws.Range("B134").Value = 3924
$ws.Range("C134").Value = 47
$ws.Range("D134").Value = 998
$ws.Range("E134").Value = 2745
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 181

# Row 148: Islandia -> Guyana
$ws.Range("A148").Value = "Guyana"
$ws.Range("B148").Value = 2535
$ws.Range("C148").Value = 98
$ws.Range("D148").Value = 1464
$ws.Range("E148").Value = 1002
$ws.Range("G148").Value = 2
$ws.Range("H148").Value = 69

# Row 149: Guyana -> Islandia
$ws.Range("A149").Value = "Islandia"
$ws.Range("B149").Value = 2476
$ws.Range("C149").Value = 57
$ws.Range("D149").Value = 2142
$ws.Range("E149").Value = 324
$ws.Range("H149").Value = 10

# Row 180: Mauricio -> San Martin (Parte Francesa)
$ws.Range("A180").Value = "San Martin (Parte Francesa)"
$ws.Range("C180").Value = 37
$ws.Range("D180").Value = 273
$ws.Range("E180").Value = 86
$ws.Range("G180").Value = 2
$ws.Range("H180").Value = 8

# Row 181: Eritrea -> Mauricio
$ws.Range("A181").Value = "Mauricio"
$ws.Range("B181").Value = 367
$ws.Range("D181").Value = 339
$ws.Range("E181").Value = 18
$ws.Range("H181").Value = 10

# Row 182: Gibraltar -> Eritrea
$ws.Range("A182").Value = "Eritrea"
$ws.Range("B182").Value = 364
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 309
$ws.Range("E182").Value = 55

# Row 183: Isla de Man -> Gibraltar
$ws.Range("A183").Value = "Gibraltar"
$ws.Range("B183").Value = 357
$ws.Range("C183").Value = 2
$ws.Range("D183").Value = 324
$ws.Range("E183").Value = 33
$ws.Range("H183").Value = 0

# Row 184: San Martin (Parte Francesa) -> Isla de Man
$ws.Range("A184").Value = "Isla de Man"
$ws.Range("B184").Value = 340
$ws.Range("D184").Value = 312
$ws.Range("E184").Value = 4
$ws.Range("H184").Value = 24

# Row 201: Puerto Rico -> San Bartolome
$ws.Range("A201").Value = "San Bartolome"
$ws.Range("B201").Value = 45
$ws.Range("C201").Value = 22
$ws.Range("D201").Value = 25
$ws.Range("E201").Value = 20
$ws.Range("H201").Value = 0

# Row 202: Guam -> Puerto Rico
$ws.Range("A202").Value = "Puerto Rico"
$ws.Range("B202").Value = 39
$ws.Range("D202").Value = 1
$ws.Range("E202").Value = 36
$ws.Range("H202").Value = 2

# Row 203: Fiyi -> Guam
$ws.Range("A203").Value = "Guam"
$ws.Range("D203").Value = 0
$ws.Range("E203").Value = 31
$ws.Range("H203").Value = 1

# Row 204: Santa Lucia -> Fiyi
$ws.Range("A204").Value = "Fiyi"
$ws.Range("B204").Value = 32
$ws.Range("D204").Value = 28
$ws.Range("E204").Value = 2
$ws.Range("H204").Value = 2

# Row 205: Timor Oriental -> Santa Lucia
$ws.Range("A205").Value = "Santa Lucia"
$ws.Range("D205").Value = 26
$ws.Range("E205").Value = 1

# Row 206: Nueva Caledonia -> Timor Oriental
$ws.Range("A206").Value = "Timor Oriental"
$ws.Range("B206").Value = 27
$ws.Range("D206").Value = 27

# Row 207: Dominica -> Nueva Caledonia
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("B207").Value = 26
$ws.Range("D207").Value = 26
$ws.Range("E207").Value = 0

# Row 208: Granada -> Dominica
$ws.Range("A208").Value = "Dominica"
$ws.Range("D208").Value = 18
$ws.Range("E208").Value = 6

# Row 209: San Bartolome -> Granada
$ws.Range("A209").Value = "Granada"
$ws.Range("B209").Value = 24
$ws.Range("D209").Value = 24
$ws.Range("E209").Value = 0

# Row 214: Montserrat -> Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215: Islas Malvinas -> Montserrat
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
